# Generate Report for Handoff
#
# The localization-status report is regenerated: every language sheet
# (and the Overview rollup) that was showing the old "Handed back: in
# sync with en-US" status now reads "Ready for handoff", the associated
# timestamps move forward a few seconds, and the Status/summary columns
# are narrowed to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-16 16:54:36"

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-16 16:54:31"

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-16 16:54:36"

# --- Re-fit the status columns now that the text is shorter ---------
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
